# Financials update: insert two new quarterly columns (new D, new E) before the
# existing data block, shifting the old D:K columns to F:M, then populate the
# two new columns with the new quarters' figures and correct two mis-keyed
# historical figures (rows 89 and 102) that moved to columns H/I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns at D (old D:K -> new F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Carry the number-format styling from column F (the old column D, now
#    shifted right) into the two freshly inserted columns for every row that
#    participates in the data block. Done per contiguous block (rows 36 and
#    78 have no row entries at all in the sheet, so they're deliberately
#    skipped to avoid manufacturing stray cells there).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the two new columns (D = newest quarter, E = next quarter) with
#    the values for every data row. Rows holding "NA" get the text "NA";
#    everything else gets its numeric value.
$newData = @"
7,43465,43373
8,9600,9300
9,5300,4500
10,4300,4800
12,NA,NA
13,0,0
14,0,0
15,0,0
17,8600,5600
18,1000,3700
20,-8900,-2200
21,NA,NA
22,0,0
23,-7900,1500
24,0,0
25,0,0
26,-7900,1500
27,-7900,1500
28,0,0
29,0,0
30,0,0
31,0,0
32,8900,2200
33,-7900,1500
34,0,0
35,-7900,1500
38,43465,43373
41,30100,29400
42,0,0
43,2600,3500
44,0,0
45,0,0
46,0,0
47,454000,383000
48,0,0
49,0,0
50,0,0
51,0,0
52,15600,25700
53,0,0
54,504000,442800
57,0,0
58,0,0
59,27600,18800
60,0,0
61,307500,243500
62,0,0
63,0,0
64,0,0
65,0,0
66,335100,262200
68,0,0
69,0,0
70,0,0
71,0,0
72,-80200,-61300
73,0,0
74,0,0
75,0,0
76,168900,180500
77,0,0
80,43465,43373
81,-7900,1500
83,0,0
84,0,0
85,0,0
86,0,0
87,0,0
88,0,0
89,-66900,19200
91,0,0
92,0,0
93,0,0
94,0,0
96,-3700,-4500
97,0,0
98,0,0
99,0,0
100,57600,2000
101,0,0
102,-9300,21200
"@

$lines = $newData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $rowNum = [int]$parts[0]
    $dText = $parts[1].Trim()
    $eText = $parts[2].Trim()

    if ($dText -eq "NA") {
        $ws.Cells.Item($rowNum, 4).Value2 = "NA"
    } else {
        $ws.Cells.Item($rowNum, 4).Value2 = [double]$dText
    }

    if ($eText -eq "NA") {
        $ws.Cells.Item($rowNum, 5).Value2 = "NA"
    } else {
        $ws.Cells.Item($rowNum, 5).Value2 = [double]$eText
    }
}

# 4. Two historical figures were re-keyed during this update. After the
#    shift they live in columns H (old F) and I (old G) of rows 89 and 102.
$ws.Cells.Item(89, 8).Value2 = -36500
$ws.Cells.Item(89, 9).Value2 = -800
$ws.Cells.Item(102, 8).Value2 = -9800
$ws.Cells.Item(102, 9).Value2 = -6300

Write-Output "Financials updated"
